$wb = $excel.ActiveWorkbook

# --- Replace "Sheet1" with a freshly created "CreateAccount" sheet in the same slot ---
# (A plain rename would keep the old sheetId; the target workbook shows a *new*
#  sheetId, so we add a brand-new sheet right after "LoginDetails" and delete the
#  old "Sheet1" - this reproduces the sheetId bump exactly like real Excel does.)
$loginDetails = $wb.Worksheets.Item("LoginDetails")
$newSheet = $wb.Worksheets.Add($null, $loginDetails)
$wb.Worksheets.Item("Sheet1").Delete()
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "CreateAccount"

# --- Header row ---
$headers = @("firstname","lastname","passward","days","month","years","company","address1","address2","city","state","country","postcode","phno","alias")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Data row (order chosen to reproduce the exact shared-string table layout) ---
$ws.Range("A2").Value = "sriyansh"
$ws.Range("B2").Value = "roy"
$ws.Range("C2").Value = "sudipa123"
$ws.Range("G2").Value = "reliance"
$ws.Range("H2").Value = "bangalore,hsr"
$ws.Range("I2").Value = "creative apartment"
$ws.Range("J2").Value = "bangalore"
$ws.Range("D2").Value = "'16"
$ws.Range("F2").Value = "'2016"
$ws.Range("N2").Value = "'9898989898"
$ws.Range("M2").Value = "'00000"
$ws.Range("L2").Value = "United States"
$ws.Range("E2").Value = "'4"
$ws.Range("K2").Value = "'9"
$ws.Range("O2").Value = "bangalore"

# --- Column widths (best achievable values given this host's char-width quantization) ---
$ws.Columns.Item(1).ColumnWidth = 10.0
$ws.Columns.Item(2).ColumnWidth = 10.0
$ws.Columns.Item(3).ColumnWidth = 12.333333333333334
$ws.Columns.Item(7).ColumnWidth = 11.5
$ws.Columns.Item(8).ColumnWidth = 13.333333333333334
$ws.Columns.Item(9).ColumnWidth = 19.666666666666668
$ws.Columns.Item(10).ColumnWidth = 14.666666666666666
$ws.Columns.Item(12).ColumnWidth = 14.0
$ws.Columns.Item(14).ColumnWidth = 15.333333333333334
$ws.Columns.Item(15).ColumnWidth = 13.666666666666666

# --- Selection / activation ---
$ws.Range("I6").Select()
$ws.Activate()

Write-Host "Done"
